$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 7 with the two missing trailing columns
$ws.Cells.Item(7, 24).Value = 0.42000000000000171
$ws.Cells.Item(7, 25).Value = "Up"

# Add new row 8 of data
$ws.Cells.Item(8, 1).Value = 42649.879930555559
$ws.Cells.Item(8, 1).NumberFormat = "m/d/yyyy h:mm"
$ws.Cells.Item(8, 2).Value = -3
$ws.Cells.Item(8, 3).Value = "Neutral"
$ws.Cells.Item(8, 4).Value = 6
$ws.Cells.Item(8, 5).Value = 5975
$ws.Cells.Item(8, 6).Value = 681
$ws.Cells.Item(8, 7).Value = 50
$ws.Cells.Item(8, 8).Value = 48
$ws.Cells.Item(8, 9).Value = 65
$ws.Cells.Item(8, 10).Value = 34
$ws.Cells.Item(8, 11).Value = 7497
$ws.Cells.Item(8, 12).Value = 82
$ws.Cells.Item(8, 13).Value = 80
$ws.Cells.Item(8, 14).Value = 29
$ws.Cells.Item(8, 15).Value = 15
$ws.Cells.Item(8, 16).Value = "Bag"
$ws.Cells.Item(8, 17).Value = 47.321424984051369
$ws.Cells.Item(8, 18).Value = 0
$ws.Cells.Item(8, 19).Value = 0.083400000000000002
$ws.Cells.Item(8, 19).NumberFormat = "0.00%"
$ws.Cells.Item(8, 20).Value = -0.0061999999999999998
$ws.Cells.Item(8, 20).NumberFormat = "0.00%"
$ws.Cells.Item(8, 21).Value = 2.31
$ws.Cells.Item(8, 22).Value = "N/A"
$ws.Cells.Item(8, 23).Value = 0

# Column C got a touch wider to fit the new "Neutral" text
$ws.Columns.Item(3).ColumnWidth = 5.6
